$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 211, shifting existing rows 211..293 down to 212..294
$ws.Rows("211:211").Insert()

# Populate the newly inserted row 211 with the new data record.
# Columns A,B,C,E,F,G,H,I,J,K,L,Q,R,T keep the same values as the (now shifted) row 212,
# while D, M, N, O, P, S hold new values.
$ws.Range("A211").Value = 5
$ws.Range("B211").Value = "Macroferia Regional de Talca"
$ws.Range("C211").Value = "Maule"
$ws.Range("D211").Value = 44837
$ws.Range("D211").NumberFormat = $ws.Range("D212").NumberFormat
$ws.Range("E211").Value = 7
$ws.Range("F211").Value = "Fruta"
$ws.Range("G211").Value = 100108
$ws.Range("H211").Value = "Tropicales y subtropicales"
$ws.Range("I211").Value = 100108005
$ws.Range("J211").Value = "Piña"
$ws.Range("K211").Value = "Caramelo"
$ws.Range("L211").Value = "Segunda"
$ws.Range("M211").Value = 540
$ws.Range("N211").Value = 21000
$ws.Range("O211").Value = 21000
$ws.Range("P211").Value = 21000
$ws.Range("Q211").Value = "$/caja 14 unidades"
$ws.Range("R211").Value = "Ecuador"
$ws.Range("S211").Value = 1500
$ws.Range("T211").Value = 14
